# Updates Adam15-Itgb1 LR-pair metrics with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 63.81813833333333
$ws.Cells.Item(2, 8).Value = 191.454415
$ws.Cells.Item(2, 9).Value = 0.5585681932726833
$ws.Cells.Item(2, 10).Value = 0.5585681932726834
$ws.Cells.Item(2, 13).Value = 159.4836373333333
$ws.Cells.Item(2, 14).Value = 478.450912
$ws.Cells.Item(2, 15).Value = 0.2983285084902258
$ws.Cells.Item(2, 16).Value = 0.2983285084902258
$ws.Cells.Item(2, 17).Value = 10177.94882924183
$ws.Cells.Item(2, 18).Value = 91601.53946317648
$ws.Cells.Item(2, 19).Value = 0.1666368159891198
$ws.Cells.Item(2, 20).Value = 0.1666368159891198

# Row 3
$ws.Cells.Item(3, 7).Value = 63.81813833333333
$ws.Cells.Item(3, 8).Value = 191.454415
$ws.Cells.Item(3, 9).Value = 0.5585681932726833
$ws.Cells.Item(3, 10).Value = 0.5585681932726834
$ws.Cells.Item(3, 13).Value = 172.558497
$ws.Cells.Item(3, 14).Value = 517.675491
$ws.Cells.Item(3, 15).Value = 0.3227862111630279
$ws.Cells.Item(3, 16).Value = 0.3227862111630279
$ws.Cells.Item(3, 17).Value = 11012.36203213808
$ws.Cells.Item(3, 18).Value = 99111.25828924275
$ws.Cells.Item(3, 19).Value = 0.1802981107826673
$ws.Cells.Item(3, 20).Value = 0.1802981107826674

# Row 4
$ws.Cells.Item(4, 7).Value = 63.81813833333333
$ws.Cells.Item(4, 8).Value = 191.454415
$ws.Cells.Item(4, 9).Value = 0.5585681932726833
$ws.Cells.Item(4, 10).Value = 0.5585681932726834
$ws.Cells.Item(4, 13).Value = 74.38770566666666
$ws.Cells.Item(4, 14).Value = 223.163117
$ws.Cells.Item(4, 15).Value = 0.1391489036280481
$ws.Cells.Item(4, 16).Value = 0.1391489036280482
$ws.Cells.Item(4, 17).Value = 4747.284890534616
$ws.Cells.Item(4, 18).Value = 42725.56401481155
$ws.Cells.Item(4, 19).Value = 0.07772415169539357
$ws.Cells.Item(4, 20).Value = 0.0777241516953936

# Row 5
$ws.Cells.Item(5, 7).Value = 63.81813833333333
$ws.Cells.Item(5, 8).Value = 191.454415
$ws.Cells.Item(5, 9).Value = 0.5585681932726833
$ws.Cells.Item(5, 10).Value = 0.5585681932726834
$ws.Cells.Item(5, 13).Value = 58.41461433333333
$ws.Cells.Item(5, 14).Value = 175.243843
$ws.Cells.Item(5, 15).Value = 0.1092697975759847
$ws.Cells.Item(5, 16).Value = 0.1092697975759848
$ws.Cells.Item(5, 17).Value = 3727.911938212982
$ws.Cells.Item(5, 18).Value = 33551.20744391684
$ws.Cells.Item(5, 19).Value = 0.06103463341128963
$ws.Cells.Item(5, 20).Value = 0.06103463341128965

# Row 6
$ws.Cells.Item(6, 7).Value = 63.81813833333333
$ws.Cells.Item(6, 8).Value = 191.454415
$ws.Cells.Item(6, 9).Value = 0.5585681932726833
$ws.Cells.Item(6, 10).Value = 0.5585681932726834
$ws.Cells.Item(6, 13).Value = 69.746216
$ws.Cells.Item(6, 14).Value = 209.238648
$ws.Cells.Item(6, 15).Value = 0.1304665791427133
$ws.Cells.Item(6, 16).Value = 0.1304665791427133
$ws.Cells.Item(6, 17).Value = 4451.073660914547
$ws.Cells.Item(6, 18).Value = 40059.66294823092
$ws.Cells.Item(6, 19).Value = 0.07287448139421292
$ws.Cells.Item(6, 20).Value = 0.07287448139421294

# Row 7
$ws.Cells.Item(7, 7).Value = 14.694925
$ws.Cells.Item(7, 8).Value = 44.084775
$ws.Cells.Item(7, 9).Value = 0.1286173166734377
$ws.Cells.Item(7, 10).Value = 0.1286173166734377
$ws.Cells.Item(7, 13).Value = 159.4836373333333
$ws.Cells.Item(7, 14).Value = 478.450912
$ws.Cells.Item(7, 15).Value = 0.2983285084902258
$ws.Cells.Item(7, 16).Value = 0.2983285084902258
$ws.Cells.Item(7, 17).Value = 2343.600089340533
$ws.Cells.Item(7, 18).Value = 21092.4008040648
$ws.Cells.Item(7, 19).Value = 0.03837021224920171
$ws.Cells.Item(7, 20).Value = 0.03837021224920172

# Row 8
$ws.Cells.Item(8, 7).Value = 14.694925
$ws.Cells.Item(8, 8).Value = 44.084775
$ws.Cells.Item(8, 9).Value = 0.1286173166734377
$ws.Cells.Item(8, 10).Value = 0.1286173166734377
$ws.Cells.Item(8, 13).Value = 172.558497
$ws.Cells.Item(8, 14).Value = 517.675491
$ws.Cells.Item(8, 15).Value = 0.3227862111630279
$ws.Cells.Item(8, 16).Value = 0.3227862111630279
$ws.Cells.Item(8, 17).Value = 2535.734171527725
$ws.Cells.Item(8, 18).Value = 22821.60754374952
$ws.Cells.Item(8, 19).Value = 0.04151589633897428
$ws.Cells.Item(8, 20).Value = 0.0415158963389743

# Row 9
$ws.Cells.Item(9, 7).Value = 14.694925
$ws.Cells.Item(9, 8).Value = 44.084775
$ws.Cells.Item(9, 9).Value = 0.1286173166734377
$ws.Cells.Item(9, 10).Value = 0.1286173166734377
$ws.Cells.Item(9, 13).Value = 74.38770566666666
$ws.Cells.Item(9, 14).Value = 223.163117
$ws.Cells.Item(9, 15).Value = 0.1391489036280481
$ws.Cells.Item(9, 16).Value = 0.1391489036280482
$ws.Cells.Item(9, 17).Value = 1093.121755693742
$ws.Cells.Item(9, 18).Value = 9838.095801243675
$ws.Cells.Item(9, 19).Value = 0.01789695860269033
$ws.Cells.Item(9, 20).Value = 0.01789695860269034

# Row 10
$ws.Cells.Item(10, 7).Value = 14.694925
$ws.Cells.Item(10, 8).Value = 44.084775
$ws.Cells.Item(10, 9).Value = 0.1286173166734377
$ws.Cells.Item(10, 10).Value = 0.1286173166734377
$ws.Cells.Item(10, 13).Value = 58.41461433333333
$ws.Cells.Item(10, 14).Value = 175.243843
$ws.Cells.Item(10, 15).Value = 0.1092697975759847
$ws.Cells.Item(10, 16).Value = 0.1092697975759848
$ws.Cells.Item(10, 17).Value = 858.3983765322583
$ws.Cells.Item(10, 18).Value = 7725.585388790325
$ws.Cells.Item(10, 19).Value = 0.01405398815767286
$ws.Cells.Item(10, 20).Value = 0.01405398815767287

# Row 11
$ws.Cells.Item(11, 7).Value = 14.694925
$ws.Cells.Item(11, 8).Value = 44.084775
$ws.Cells.Item(11, 9).Value = 0.1286173166734377
$ws.Cells.Item(11, 10).Value = 0.1286173166734377
$ws.Cells.Item(11, 13).Value = 69.746216
$ws.Cells.Item(11, 14).Value = 209.238648
$ws.Cells.Item(11, 15).Value = 0.1304665791427133
$ws.Cells.Item(11, 16).Value = 0.1304665791427133
$ws.Cells.Item(11, 17).Value = 1024.9154131538
$ws.Cells.Item(11, 18).Value = 9224.2387183842
$ws.Cells.Item(11, 19).Value = 0.01678026132489848
$ws.Cells.Item(11, 20).Value = 0.01678026132489848

# Row 12
$ws.Cells.Item(12, 7).Value = 17.02115633333333
$ws.Cells.Item(12, 8).Value = 51.063469
$ws.Cells.Item(12, 9).Value = 0.1489776541406249
$ws.Cells.Item(12, 10).Value = 0.1489776541406249
$ws.Cells.Item(12, 13).Value = 159.4836373333333
$ws.Cells.Item(12, 14).Value = 478.450912
$ws.Cells.Item(12, 15).Value = 0.2983285084902258
$ws.Cells.Item(12, 16).Value = 0.2983285084902258
$ws.Cells.Item(12, 17).Value = 2714.595923659303
$ws.Cells.Item(12, 18).Value = 24431.36331293373
$ws.Cells.Item(12, 19).Value = 0.04444428135814534
$ws.Cells.Item(12, 20).Value = 0.04444428135814535

# Row 13
$ws.Cells.Item(13, 7).Value = 17.02115633333333
$ws.Cells.Item(13, 8).Value = 51.063469
$ws.Cells.Item(13, 9).Value = 0.1489776541406249
$ws.Cells.Item(13, 10).Value = 0.1489776541406249
$ws.Cells.Item(13, 13).Value = 172.558497
$ws.Cells.Item(13, 14).Value = 517.675491
$ws.Cells.Item(13, 15).Value = 0.3227862111630279
$ws.Cells.Item(13, 16).Value = 0.3227862111630279
$ws.Cells.Item(13, 17).Value = 2937.145154082031
$ws.Cells.Item(13, 18).Value = 26434.30638673828
$ws.Cells.Item(13, 19).Value = 0.0480879325280083
$ws.Cells.Item(13, 20).Value = 0.0480879325280083

# Row 14
$ws.Cells.Item(14, 7).Value = 17.02115633333333
$ws.Cells.Item(14, 8).Value = 51.063469
$ws.Cells.Item(14, 9).Value = 0.1489776541406249
$ws.Cells.Item(14, 10).Value = 0.1489776541406249
$ws.Cells.Item(14, 13).Value = 74.38770566666666
$ws.Cells.Item(14, 14).Value = 223.163117
$ws.Cells.Item(14, 15).Value = 0.1391489036280481
$ws.Cells.Item(14, 16).Value = 0.1391489036280482
$ws.Cells.Item(14, 17).Value = 1266.164767430319
$ws.Cells.Item(14, 18).Value = 11395.48290687287
$ws.Cells.Item(14, 19).Value = 0.0207300772387465
$ws.Cells.Item(14, 20).Value = 0.02073007723874651

# Row 15
$ws.Cells.Item(15, 7).Value = 17.02115633333333
$ws.Cells.Item(15, 8).Value = 51.063469
$ws.Cells.Item(15, 9).Value = 0.1489776541406249
$ws.Cells.Item(15, 10).Value = 0.1489776541406249
$ws.Cells.Item(15, 13).Value = 58.41461433333333
$ws.Cells.Item(15, 14).Value = 175.243843
$ws.Cells.Item(15, 15).Value = 0.1092697975759847
$ws.Cells.Item(15, 16).Value = 0.1092697975759848
$ws.Cells.Item(15, 17).Value = 994.2842827190408
$ws.Cells.Item(15, 18).Value = 8948.558544471367
$ws.Cells.Item(15, 19).Value = 0.01627875811129115
$ws.Cells.Item(15, 20).Value = 0.01627875811129115

# Row 16
$ws.Cells.Item(16, 7).Value = 17.02115633333333
$ws.Cells.Item(16, 8).Value = 51.063469
$ws.Cells.Item(16, 9).Value = 0.1489776541406249
$ws.Cells.Item(16, 10).Value = 0.1489776541406249
$ws.Cells.Item(16, 13).Value = 69.746216
$ws.Cells.Item(16, 14).Value = 209.238648
$ws.Cells.Item(16, 15).Value = 0.1304665791427133
$ws.Cells.Item(16, 16).Value = 0.1304665791427133
$ws.Cells.Item(16, 17).Value = 1187.161246194435
$ws.Cells.Item(16, 18).Value = 10684.45121574991
$ws.Cells.Item(16, 19).Value = 0.01943660490443361
$ws.Cells.Item(16, 20).Value = 0.01943660490443362

# Row 17
$ws.Cells.Item(17, 7).Value = 0.7288956666666667
$ws.Cells.Item(17, 8).Value = 2.186687
$ws.Cells.Item(17, 9).Value = 0.006379658608775693
$ws.Cells.Item(17, 10).Value = 0.006379658608775693
$ws.Cells.Item(17, 13).Value = 159.4836373333333
$ws.Cells.Item(17, 14).Value = 478.450912
$ws.Cells.Item(17, 15).Value = 0.2983285084902258
$ws.Cells.Item(17, 16).Value = 0.2983285084902258
$ws.Cells.Item(17, 17).Value = 116.2469321565049
$ws.Cells.Item(17, 18).Value = 1046.222389408544
$ws.Cells.Item(17, 19).Value = 0.001903234037432881
$ws.Cells.Item(17, 20).Value = 0.001903234037432881

# Row 18
$ws.Cells.Item(18, 7).Value = 0.7288956666666667
$ws.Cells.Item(18, 8).Value = 2.186687
$ws.Cells.Item(18, 9).Value = 0.006379658608775693
$ws.Cells.Item(18, 10).Value = 0.006379658608775693
$ws.Cells.Item(18, 13).Value = 172.558497
$ws.Cells.Item(18, 14).Value = 517.675491
$ws.Cells.Item(18, 15).Value = 0.3227862111630279
$ws.Cells.Item(18, 16).Value = 0.3227862111630279
$ws.Cells.Item(18, 17).Value = 125.777140709813
$ws.Cells.Item(18, 18).Value = 1131.994266388317
$ws.Cells.Item(18, 19).Value = 0.0020592658308403
$ws.Cells.Item(18, 20).Value = 0.0020592658308403

# Row 19
$ws.Cells.Item(19, 7).Value = 0.7288956666666667
$ws.Cells.Item(19, 8).Value = 2.186687
$ws.Cells.Item(19, 9).Value = 0.006379658608775693
$ws.Cells.Item(19, 10).Value = 0.006379658608775693
$ws.Cells.Item(19, 13).Value = 74.38770566666666
$ws.Cells.Item(19, 14).Value = 223.163117
$ws.Cells.Item(19, 15).Value = 0.1391489036280481
$ws.Cells.Item(19, 16).Value = 0.1391489036280482
$ws.Cells.Item(19, 17).Value = 54.22087631370878
$ws.Cells.Item(19, 18).Value = 487.987886823379
$ws.Cells.Item(19, 19).Value = 0.0008877225009323765
$ws.Cells.Item(19, 20).Value = 0.0008877225009323767

# Row 20
$ws.Cells.Item(20, 7).Value = 0.7288956666666667
$ws.Cells.Item(20, 8).Value = 2.186687
$ws.Cells.Item(20, 9).Value = 0.006379658608775693
$ws.Cells.Item(20, 10).Value = 0.006379658608775693
$ws.Cells.Item(20, 13).Value = 58.41461433333333
$ws.Cells.Item(20, 14).Value = 175.243843
$ws.Cells.Item(20, 15).Value = 0.1092697975759847
$ws.Cells.Item(20, 16).Value = 0.1092697975759848
$ws.Cells.Item(20, 17).Value = 42.57815925757122
$ws.Cells.Item(20, 18).Value = 383.203433318141
$ws.Cells.Item(20, 19).Value = 0.0006971040047848084
$ws.Cells.Item(20, 20).Value = 0.0006971040047848085

# Row 21
$ws.Cells.Item(21, 7).Value = 0.7288956666666667
$ws.Cells.Item(21, 8).Value = 2.186687
$ws.Cells.Item(21, 9).Value = 0.006379658608775693
$ws.Cells.Item(21, 10).Value = 0.006379658608775693
$ws.Cells.Item(21, 13).Value = 69.746216
$ws.Cells.Item(21, 14).Value = 209.238648
$ws.Cells.Item(21, 15).Value = 0.1304665791427133
$ws.Cells.Item(21, 16).Value = 0.1304665791427133
$ws.Cells.Item(21, 17).Value = 50.83771460879734
$ws.Cells.Item(21, 18).Value = 457.539431479176
$ws.Cells.Item(21, 19).Value = 0.0008323322347853262
$ws.Cells.Item(21, 20).Value = 0.0008323322347853265

# Row 22
$ws.Cells.Item(22, 7).Value = 17.989968
$ws.Cells.Item(22, 8).Value = 53.969904
$ws.Cells.Item(22, 9).Value = 0.1574571773044783
$ws.Cells.Item(22, 10).Value = 0.1574571773044783
$ws.Cells.Item(22, 13).Value = 159.4836373333333
$ws.Cells.Item(22, 14).Value = 478.450912
$ws.Cells.Item(22, 15).Value = 0.2983285084902258
$ws.Cells.Item(22, 16).Value = 0.2983285084902258
$ws.Cells.Item(22, 17).Value = 2869.105532150273
$ws.Cells.Item(22, 18).Value = 25821.94978935245
$ws.Cells.Item(22, 19).Value = 0.04697396485632604
$ws.Cells.Item(22, 20).Value = 0.04697396485632604

# Row 23
$ws.Cells.Item(23, 7).Value = 17.989968
$ws.Cells.Item(23, 8).Value = 53.969904
$ws.Cells.Item(23, 9).Value = 0.1574571773044783
$ws.Cells.Item(23, 10).Value = 0.1574571773044783
$ws.Cells.Item(23, 13).Value = 172.558497
$ws.Cells.Item(23, 14).Value = 517.675491
$ws.Cells.Item(23, 15).Value = 0.3227862111630279
$ws.Cells.Item(23, 16).Value = 0.3227862111630279
$ws.Cells.Item(23, 17).Value = 3104.321839158096
$ws.Cells.Item(23, 18).Value = 27938.89655242286
$ws.Cells.Item(23, 19).Value = 0.05082500568253765
$ws.Cells.Item(23, 20).Value = 0.05082500568253765

# Row 24
$ws.Cells.Item(24, 7).Value = 17.989968
$ws.Cells.Item(24, 8).Value = 53.969904
$ws.Cells.Item(24, 9).Value = 0.1574571773044783
$ws.Cells.Item(24, 10).Value = 0.1574571773044783
$ws.Cells.Item(24, 13).Value = 74.38770566666666
$ws.Cells.Item(24, 14).Value = 223.163117
$ws.Cells.Item(24, 15).Value = 0.1391489036280481
$ws.Cells.Item(24, 16).Value = 0.1391489036280482
$ws.Cells.Item(24, 17).Value = 1338.232444536752
$ws.Cells.Item(24, 18).Value = 12044.09200083077
$ws.Cells.Item(24, 19).Value = 0.02190999359028534
$ws.Cells.Item(24, 20).Value = 0.02190999359028534

# Row 25
$ws.Cells.Item(25, 7).Value = 17.989968
$ws.Cells.Item(25, 8).Value = 53.969904
$ws.Cells.Item(25, 9).Value = 0.1574571773044783
$ws.Cells.Item(25, 10).Value = 0.1574571773044783
$ws.Cells.Item(25, 13).Value = 58.41461433333333
$ws.Cells.Item(25, 14).Value = 175.243843
$ws.Cells.Item(25, 15).Value = 0.1092697975759847
$ws.Cells.Item(25, 16).Value = 0.1092697975759848
$ws.Cells.Item(25, 17).Value = 1050.877042589008
$ws.Cells.Item(25, 18).Value = 9457.893383301072
$ws.Cells.Item(25, 19).Value = 0.01720531389094628
$ws.Cells.Item(25, 20).Value = 0.01720531389094628

# Row 26
$ws.Cells.Item(26, 7).Value = 17.989968
$ws.Cells.Item(26, 8).Value = 53.969904
$ws.Cells.Item(26, 9).Value = 0.1574571773044783
$ws.Cells.Item(26, 10).Value = 0.1574571773044783
$ws.Cells.Item(26, 13).Value = 69.746216
$ws.Cells.Item(26, 14).Value = 209.238648
$ws.Cells.Item(26, 15).Value = 0.1304665791427133
$ws.Cells.Item(26, 16).Value = 0.1304665791427133
$ws.Cells.Item(26, 17).Value = 1254.732193961088
$ws.Cells.Item(26, 18).Value = 11292.58974564979
$ws.Cells.Item(26, 19).Value = 0.02054289928438296
$ws.Cells.Item(26, 20).Value = 0.02054289928438296
